# Start version 5 for EU data
# Adds a new "Group1" worksheet (after "NO") with wave-1 rows for AT, DK, ES,
# FR, IT, PL, PT, makes it the active sheet/tab, and tweaks a couple of
# leftover view-selection bits on the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the BE sheet's lingering cell selection (view-state only).
# ---------------------------------------------------------------------------
$beSheet = $wb.Worksheets.Item("BE")
$beSheet.Range("A1:M10").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Group1" worksheet as the last tab and make it active.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Group1"

# Header row - reuses existing shared strings from the other country sheets.
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "survey_version"
$ws.Range("C1").Value = "locked"
$ws.Range("D1").Value = "week"
$ws.Range("E1").Value = "panel"
$ws.Range("F1").Value = "wave"
$ws.Range("G1").Value = "date_recieved"
$ws.Range("H1").Value = "spss_name"
$ws.Range("I1").Value = "r_name"
$ws.Range("J1").Value = "r_saved"
$ws.Range("K1").Value = "cleaned"
$ws.Range("L1").Value = "combined"
$ws.Range("M1").Value = "cleaned_by"

# Row 2 - AT (typed first: spss_name filename, then the country code, matching
# the original authoring order so new shared strings land at the same index).
$ws.Range("H2").Value = "20-030971_AT_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("A2").Value = "at"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "A"
$ws.Range("F2").Value = 1
$ws.Range("J2").Value = 1

# Rows 3-8 - country codes filled down column A first ...
$ws.Range("A3").Value = "dk"
$ws.Range("A4").Value = "es"
$ws.Range("A5").Value = "fr"
$ws.Range("A6").Value = "it"
$ws.Range("A7").Value = "pl"
$ws.Range("A8").Value = "pt"

# ... then the spss_name filenames down column H.
$ws.Range("H3").Value = "20-030971_DK_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("H4").Value = "20-030971_ES_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("H5").Value = "20-030971_FR_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("H6").Value = "20-030971_IT_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("H7").Value = "20-030971_PL_Wave1_Final_v1_110121_IntClientuse"
$ws.Range("H8").Value = "20-030971_PT_Wave1_Final_v1_110121_IntClientuse"

# Remaining numeric/text columns for rows 3-8 (same values as row 2).
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("B$r").Value = 5
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = "A"
    $ws.Range("F$r").Value = 1
    $ws.Range("J$r").Value = 1
}

# Date received column (2021-01-11 for every row) - pasted from an existing
# date cell so the short-date style is reused instead of creating a new one.
$ws.Range("G2:G8").Value = 44207
$wb.Worksheets.Item("NO").Range("G2").Copy() | Out-Null
$ws.Range("G2:G10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# r_name formula. Row 2 was typed by itself (quirky G3 reference kept as-is);
# rows 3-9 were filled as one shared-formula group anchored at I3.
$ws.Range("I2").Formula = '=A2&"_"&"wk"&TEXT(D2,"00")&"_"&YEAR(G3)&TEXT(G3,"MM")&TEXT(G3,"DD")&"_p"&E2&"_wv"&TEXT(F2,"00")&""'
$ws.Range("I3:I9").Formula = '=A3&"_"&"wk"&TEXT(D3,"00")&"_"&YEAR(G3)&TEXT(G3,"MM")&TEXT(G3,"DD")&"_p"&E3&"_wv"&TEXT(F3,"00")&""'
$ws.Range("I9").ClearContents()

# Column widths (G/H/I) to roughly match the auto-fit widths of the sibling
# country sheets.
$ws.Columns.Item(7).ColumnWidth = 9.8
$ws.Columns.Item(8).ColumnWidth = 58.15
$ws.Columns.Item(9).ColumnWidth = 26.5

$ws.Range("F9").Select()
